# "Generate Report for Archive"
#
# The localization-status report is being regenerated: the shared status
# string "Ready for handoff" becomes "In Translation" everywhere it is
# used (Overview sheet's per-language status columns, plus each
# language's own status column), and the status columns are narrowed to
# match the new (shorter) auto-fit text width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update every cell currently showing the old status text ---------

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- Narrow the affected status columns to fit the shorter text ------
# (ColumnWidth is in characters; 12.5 is the calibrated input that
# reproduces the auto-fit result for the new, shorter status text.)

$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
